# "ui for configuration of TimeFlip device"
# Fill in the real time-tracking entries for the rows that were only
# templated (empty) before, fix up the "REST API" label, grow the sheet's
# empty-row buffer by five rows, and extend the activity dropdown.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datenerfassung")

# ---------------------------------------------------------------------
# D33: "REST API + Refactoring" -> "REST API"
# ---------------------------------------------------------------------
$ws.Range("D33").Value = "REST API"

# ---------------------------------------------------------------------
# Rows 36-46: these were blank placeholder rows; fill them with the real
# logged entries (date, duration, primary activity, note).
# ---------------------------------------------------------------------
$ws.Range("A36").Value = 44298
$ws.Range("B36").Value = 0.041666666666666664
$ws.Range("C36").Value = "LV-Einheit"

$ws.Range("A37").Value = 44300
$ws.Range("B37").Value = 0.041666666666666664
$ws.Range("C37").Value = "Koordination und Projektmanagement"
$ws.Range("D37").Value = "Teammeeting"

$ws.Range("A38").Value = 44301
$ws.Range("B38").Value = 0.125
$ws.Range("C38").Value = "Implementierung"

$ws.Range("A39").Value = 44303
$ws.Range("B39").Value = 0.08333333333333333
$ws.Range("C39").Value = "Software/System Design und Architektur"
$ws.Range("D39").Value = "Besprechung zum Kommunikationsprotokoll zwischen Raspberry Pi und Backend"

$ws.Range("A40").Value = 44303
$ws.Range("B40").Value = 0.125
$ws.Range("C40").Value = "Implementierung"

$ws.Range("A41").Value = 44305
$ws.Range("B41").Value = 0.20833333333333334
$ws.Range("C41").Value = "Implementierung"
$ws.Range("D41").Value = "Zuweisung von Würfelseiten"

$ws.Range("A42").Value = 44306
$ws.Range("B42").Value = 0.041666666666666664
$ws.Range("C42").Value = "Koordination und Projektmanagement"
$ws.Range("D42").Value = "Teammeeting"

$ws.Range("A43").Value = 44307
$ws.Range("B43").Value = 0.041666666666666664
$ws.Range("C43").Value = "Implementierung"
$ws.Range("D43").Value = "Raspberry Pi"

$ws.Range("A44").Value = 44307
$ws.Range("B44").Value = 0.08333333333333333
$ws.Range("C44").Value = "Implementierung"
$ws.Range("D44").Value = "Zuweisung von Würfelseiten"

$ws.Range("A45").Value = 44310
$ws.Range("B45").Value = 0.25
$ws.Range("C45").Value = "Implementierung"
$ws.Range("D45").Value = "Zuweisung von Würfelseiten"

$ws.Range("A46").Value = 44311
$ws.Range("B46").Value = 0.041666666666666664
$ws.Range("C46").Value = "Implementierung"
$ws.Range("D46").Value = "Zuweisung von Würfelseiten"

# ---------------------------------------------------------------------
# Grow the trailing blank-row buffer by five rows (the sheet always keeps
# a pool of empty, pre-formatted rows after the last data row, terminated
# by a lone dated sentinel row). Insert five rows just above the sentinel
# so the sentinel moves from row 1011 to row 1016.
# ---------------------------------------------------------------------
$ws.Range("A1011:A1015").EntireRow.Insert()

# ---------------------------------------------------------------------
# Scroll/selection bookkeeping to match where the author was last working.
# ---------------------------------------------------------------------
$ws.Range("D47").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 15
$win.ScrollColumn = 3
